$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing text entries (sharedStrings) ---

# Append text to the "Compte rendu fait par M. Egger..." note (row 31, col B)
$ws.Range("B31").Value2 = "Compte rendu fait par M. Egger sur ma documentation de projet. Nous en avons discuté et j'ai commencé à améliorer les points qui sont sortis durant la disscusion. J'ai ensuite passé mon document sur un nouveau caneva qu'on nous a distribué."

# Row 31: the "1 périodes" entry in C31 becomes "2 périodes" (reuses existing shared string)
$ws.Range("C31").Value2 = "2 périodes"

# Row 31's height grows from 30 to 45 to fit the longer text
$ws.Rows(31).RowHeight = 45

# --- Add the new journal entry as row 32 ---

# Duplicate formatting of row 31 onto the new row 32 first
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A32").Value2 = 43161
$ws.Range("B32").Value2 = "J'ai commencé à mettre des illustrations pour chaque article quand on arrive sur la page product-details.php. Ensuite j'ai commencé à faire la requête pour l'ajout des articles dans le panier. "
$ws.Range("C32").Value2 = "2 périodes"

$ws.Rows(32).RowHeight = 45

# --- Update the active selection to reflect the next empty row ---
$ws.Range("C33").Select() | Out-Null
